$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that correspond to "2020-06" and "2020-07" (old rows 113 and 114).
# Everything below shifts up by two rows.
$ws.Rows("113:114").Delete()

# Column A used to hold text labels like "2011-01"; it now holds the actual date
# (first day of the month) as a real date serial number, formatted as a date-time.
$dateSerials = @(
    40544,40575,40603,40634,40664,40695,40725,40756,40787,40817,
    40848,40878,40909,40940,40969,41000,41030,41061,41091,41122,
    41153,41183,41214,41244,41275,41306,41334,41365,41395,41426,
    41456,41487,41518,41548,41579,41609,41640,41671,41699,41730,
    41760,41791,41821,41852,41883,41913,41944,41974,42005,42036,
    42064,42095,42125,42156,42186,42217,42248,42278,42309,42339,
    42370,42401,42430,42461,42491,42522,42552,42583,42614,42644,
    42675,42705,42736,42767,42795,42826,42856,42887,42917,42948,
    42979,43009,43040,43070,43101,43132,43160,43191,43221,43252,
    43282,43313,43344,43374,43405,43435,43466,43497,43525,43556,
    43586,43617,43647,43678,43709,43739,43770,43800,43831,43862,
    43891,44044,44075,44105,44136,44166,44197,44228,44256,44287,
    44317,44348,44378,44409,44440,44470,44501,44531,44562,44593,
    44621,44652,44682,44713,44743,44774,44805,44835,44866,44896,
    44927
)

$firstRow = 2
$lastRow = $firstRow + $dateSerials.Length - 1

for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $dateSerials[$i]
}

# Apply the number format. The first cell is formatted twice (lower-case form first,
# then the final upper-case form) which is what registers both custom number formats
# (164 = "yyyy-mm-dd h:mm:ss", 165 = "YYYY-MM-DD HH:MM:SS") while the cell itself ends
# up using the final (165) format. All remaining cells are then set straight to the
# final format, re-using the same style.
$ws.Cells.Item($firstRow, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item($firstRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
